$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '56.828.09'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '2.504.89'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '496.88'
$ws.Range("E5").Value = '  +2.90%  '
$ws.Range("D6").Value = '154.26'
$ws.Range("E6").Value = '  +9.42%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("D9").Value = '2.512.64'
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("E10").Value = '  +5.41%  '
$ws.Range("D11").Value = '0.0995'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("E12").Value = '  +2.73%  '
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").Value = '2.933.66'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '56.959.51'
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").Value = '21.58'
$ws.Range("E16").Value = '  +4.71%  '
$ws.Range("D17").Value = '0.0000138'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").Value = '2.513.86'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '4.57'
$ws.Range("E19").Value = '  +4.38%  '
$ws.Range("D20").Value = '10.37'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").Value = '325.23'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  +3.89%  '
$ws.Range("D24").Value = '59.18'
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").Value = '0.413'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").Value = '2.614.33'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("D30").Value = '0.0₃0824'
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").Value = '152.59'
$ws.Range("E32").Value = '  +2.24%  '
$ws.Range("D33").Value = '18.45'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("E34").Value = '  +3.45%  '
$ws.Range("D35").Value = '5.28'
$ws.Range("E35").Value = '  +1.77%  '
$ws.Range("E36").Value = '  +4.47%  '
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").Value = '0.880'
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("E39").Value = '  +5.97%  '
$ws.Range("D40").Value = '34.26'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").Value = '0.0567'
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("D43").Value = '0.618'
$ws.Range("E43").Value = '  +1.50%  '
$ws.Range("D44").Value = '0.994'
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '4.98'
$ws.Range("E45").Value = '  +7.58%  '
$ws.Range("D46").Value = '268.43'
$ws.Range("E46").Value = '  +6.18%  '
$ws.Range("D47").Value = '0.0933'
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").Value = '0.0232'
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").Value = '17.91'
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("D51").Value = '1.912.60'
$ws.Range("E51").Value = '  -3.28%  '
